# Commit: "pre or prod arguments"
# The edit reorders the comma-separated activity-code lists in the
# "Werkzaamheden" column (and a couple of related date/time cells) so the
# values match between the "Pre" and "Prod" status logs. No rows/columns are
# inserted or removed - only the text/values of specific cells change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value  = 'SlopenVanEenBouwwwerkOfDeelDaarvanOfHetVerwijderenVanAsbest, BouwwerkOnderhouden, BouwwerkReinigenOfConserveren'
$ws.Range("E13").Value = 'GietwaterKlaarmaken, GewassenSpoelen, DrinkwaterVoorVeeKlaarmaken, FruitSorteren'
$ws.Range("E16").Value = 'VisvijverExploiteren, BedrijfVoorTelenEnKwekenVanWaterplantenOfWaterdieren'
$ws.Range("E24").Value = 'BrugPlaatsen, BrugVeranderen, BrugWeghalen'
$ws.Range("E26").Value = 'GevaarlijkeStoffenOpslaanInVerpakking, GasOpslaanInEenOpslagtank, OntplofbareStoffenOfVoorwerpenOpslaanEtc'

$ws.Range("G30").Value = '20-03-2024 16:38:36'
$ws.Range("H30").Value = '20-03-2024 16:25:46'
$ws.Range("I30").Value = '20-03-2024 16:25:43'

$ws.Range("E33").Value = 'BeschoeiingWeghalen, BeschoeiingVeranderen, BeschoeiingPlaatsen'
$ws.Range("E35").Value = 'LozingsvoorzieningPlaatsen, OnttrekkingsvoorzieningPlaatsen'
$ws.Range("E36").Value = 'OnttrekkingsvoorzieningPlaatsen, LozingsvoorzieningWeghalen, LozingsvoorzieningPlaatsen'
$ws.Range("E40").Value = 'VerhardingAanbrengen, VerhardingWeghalen'
$ws.Range("E41").Value = 'AardgasBehandelen, ConstructiePlaatsen, BouwwerkAgrarischPlaatsen, BouwwerkInfrastructuurPlaatsen, ZonnepaneelPlaatsen, SpeeltoestelPlaatsen, TuinmeubilairPlaatsen'
$ws.Range("E42").Value = 'ConstructieWeghalen, SlopenVanEenBouwwwerkOfDeelDaarvanOfHetVerwijderenVanAsbest, ZonnepaneelPlaatsen, BouwwerkAgrarischPlaatsen, TuinmeubilairPlaatsen, SpeeltoestelPlaatsen, ConstructiePlaatsen, BouwwerkInfrastructuurPlaatsen'
$ws.Range("E43").Value = 'AanbouwPlaatsen, GebouwVeranderen, GebouwPlaatsen, AanbouwVeranderen, BouwwerkAgrarischPlaatsen, AardgasBehandelen, Windturbine'
$ws.Range("E44").Value = 'AanbouwPlaatsen, SlopenVanEenBouwwwerkOfDeelDaarvanOfHetVerwijderenVanAsbest, Windturbine, BouwwerkAgrarischPlaatsen'
$ws.Range("E47").Value = 'BoomPlanten, BoomWeghalen'
$ws.Range("E52").Value = 'PeilscheidingAanleggen, PeilscheidingVerwijderen'
$ws.Range("E54").Value = 'CivielKunstwerkPlaatsen, CivielKunstwerkWeghalen'
$ws.Range("E57").Value = 'DamPlaatsen, DamWeghalen'
$ws.Range("E63").Value = 'KabelPlaatsen, KabelWeghalen'
